# Apply the commit: add "get_support" sheet (cloned from the my_account
# template), rename "Sheet2" to "my_account", and update the relevant
# cell contents / selections.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$template = $wb.Worksheets.Item("Sheet2")

# Duplicate the "Sheet2" (my_account) sheet and place the copy before
# "Sheet1".
$template.Copy($ws1)

# Re-fetch the sheets by name: inserting the new sheet shifts indices,
# so stale index-based references must not be relied upon.
$newSheet = $wb.Worksheets.Item("Sheet2 (2)")
$newSheet.Name = "get_support"

$template = $wb.Worksheets.Item("Sheet2")
$template.Name = "my_account"

# Update the new get_support sheet with its own test-case summary data.
$newSheet.Range("C2").Value = " Test Case Summary(12-02-24)"
$newSheet.Range("D4").Value = 5
$newSheet.Range("D6").Value = 15
$newSheet.Range("D8").Value = "TC_SYM_GSF_0021"

# Update the selection remembered on the my_account sheet, then
# re-activate get_support so it stays the selected/visible tab.
$template.Range("H8").Select() | Out-Null
$newSheet.Range("D8").Select() | Out-Null
